$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 90; this shifts existing rows 90-112 down to 91-113
$ws.Rows("90:90").Insert()

# Populate the newly inserted row 90 with the new weekly record
$ws.Cells.Item(90, 1).Value = 11
$ws.Cells.Item(90, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(90, 3).Value = "Bíobío"
$ws.Cells.Item(90, 4).Value = 44642
$ws.Cells.Item(90, 5).Value = 8
$ws.Cells.Item(90, 6).Value = 100112043
$ws.Cells.Item(90, 7).Value = "Pepino ensalada"
$ws.Cells.Item(90, 8).Value = "Sin especificar"
$ws.Cells.Item(90, 9).Value = "Primera"
$ws.Cells.Item(90, 10).Value = 270
$ws.Cells.Item(90, 11).Value = 17000
$ws.Cells.Item(90, 12).Value = 18000
$ws.Cells.Item(90, 13).Value = 17556
$ws.Cells.Item(90, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(90, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(90, 16).Value = 293
$ws.Cells.Item(90, 17).Value = 60
$ws.Cells.Item(90, 18).Value = "Hortaliza"
